# C1--C2-and-C3-PowerPoint.pptx
#
# 1) Slide 16's table (3rd shape) switches to a different built-in table
#    style (tableStyleId changes from {C7615391-...} to {4DC69F61-...}).
# 2) The deck's theme is switched from the "Integral" theme to the
#    standard "Office Theme" (12 theme colours re-mapped on the slide
#    master's theme).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{4DC69F61-0429-4B2C-BC9B-CBF44F72749E}")
}

# --- 2. Swap the presentation theme from "Integral" to "Office Theme" ------
function Convert-HexToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme, in ThemeColorScheme.Item() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1..10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToOleColor $officeThemeColors[$i - 1]
}
